# Newton's method worksheet: find sqrt(28561) by iterating
#   x_{n+1} = x_n - f(x_n)/f'(x_n),  f(x) = x^2 - y,  f'(x) = 2x
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Initial guess (B3) and target value (C3), row 3 seeds the iteration
$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 28561
$ws.Range("D3").Formula = '=B3^2-$C$3'
$ws.Range("E3").Formula = '=2*B3'

# Row 4 is the first Newton step, referencing row 3
$ws.Range("B4").Formula = '=B3-(D3/E3)'
$ws.Range("D4").Formula = '=B4^2-$C$3'
$ws.Range("E4").Formula = '=2*B4'

# Rows 5-23 repeat the same pattern, each referencing the row above
for ($r = 5; $r -le 23; $r++) {
    $prev = $r - 1
    $ws.Range("B$r").Formula = "=B$prev-(D$prev/E$prev)"
    $ws.Range("D$r").Formula = "=B$r" + '^2-$C$3'
    $ws.Range("E$r").Formula = "=2*B$r"
}

# Restore the author's final selection
$ws.Range("C9").Select()
